# Refreshes the Indiana COVID deaths-by-date-by-age-group dataset tail (rows 970-1036).
# Row 970 onward is re-synced to the latest upstream snapshot: some existing counts are
# revised, the date groupings shift down by one row starting 2020-11-25 (serial 44160),
# and four brand new report dates (2020-12-01 .. 2020-12-04) are appended at the end,
# extending the sheet from A1:C1020 to A1:C1036.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is: worksheet row, date serial, age group label, covid_deaths count
$rows = @(
  @(970, 44155, "70-79", 12),
  @(971, 44155, "80+", 31),
  @(972, 44156, "40-49", 1),
  @(973, 44156, "60-69", 9),
  @(974, 44156, "70-79", 11),
  @(975, 44156, "80+", 34),
  @(976, 44157, "30-39", 1),
  @(977, 44157, "50-59", 2),
  @(978, 44157, "60-69", 8),
  @(979, 44157, "70-79", 16),
  @(980, 44157, "80+", 34),
  @(981, 44158, "60-69", 13),
  @(982, 44158, "70-79", 22),
  @(983, 44158, "80+", 34),
  @(984, 44159, "50-59", 1),
  @(985, 44159, "60-69", 9),
  @(986, 44159, "70-79", 17),
  @(987, 44159, "80+", 34),
  @(988, 44160, "30-39", 1),
  @(989, 44160, "40-49", 1),
  @(990, 44160, "50-59", 1),
  @(991, 44160, "60-69", 6),
  @(992, 44160, "70-79", 15),
  @(993, 44160, "80+", 31),
  @(994, 44161, "40-49", 3),
  @(995, 44161, "50-59", 4),
  @(996, 44161, "60-69", 9),
  @(997, 44161, "70-79", 14),
  @(998, 44161, "80+", 37),
  @(999, 44162, "50-59", 4),
  @(1000, 44162, "60-69", 8),
  @(1001, 44162, "70-79", 29),
  @(1002, 44162, "80+", 27),
  @(1003, 44163, "30-39", 2),
  @(1004, 44163, "40-49", 1),
  @(1005, 44163, "50-59", 2),
  @(1006, 44163, "60-69", 11),
  @(1007, 44163, "70-79", 11),
  @(1008, 44163, "80+", 26),
  @(1009, 44164, "50-59", 2),
  @(1010, 44164, "60-69", 7),
  @(1011, 44164, "70-79", 18),
  @(1012, 44164, "80+", 36),
  @(1013, 44165, "0-19", 1),
  @(1014, 44165, "50-59", 2),
  @(1015, 44165, "60-69", 6),
  @(1016, 44165, "70-79", 19),
  @(1017, 44165, "80+", 36),
  @(1018, 44166, "0-19", 1),
  @(1019, 44166, "40-49", 1),
  @(1020, 44166, "50-59", 2),
  @(1021, 44166, "60-69", 5),
  @(1022, 44166, "70-79", 20),
  @(1023, 44166, "80+", 32),
  @(1024, 44167, "40-49", 2),
  @(1025, 44167, "50-59", 4),
  @(1026, 44167, "60-69", 6),
  @(1027, 44167, "70-79", 15),
  @(1028, 44167, "80+", 40),
  @(1029, 44168, "50-59", 1),
  @(1030, 44168, "60-69", 6),
  @(1031, 44168, "70-79", 14),
  @(1032, 44168, "80+", 31),
  @(1033, 44169, "50-59", 1),
  @(1034, 44169, "60-69", 5),
  @(1035, 44169, "70-79", 6),
  @(1036, 44169, "80+", 21)
)

foreach ($row in $rows) {
    $r = $row[0]
    $dateSerial = $row[1]
    $ageGroup = $row[2]
    $deaths = $row[3]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $dateSerial
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $ageGroup
    $ws.Cells.Item($r, 3).Value = $deaths
}
